$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Retailer OrderPlacement")
$ws.Activate()
$ws.Columns("N").Insert()
Write-Host "done"
